$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.296.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.49"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3265"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.84"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.73%  "

$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07361"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.836"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.747"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.562.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001074"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06635"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.24%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.392"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.301.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.552"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.927"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.740.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.075"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.886"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.885"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.354"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08213"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02366"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06278"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.284"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2151"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.241"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.03"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6057"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.739"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5914"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.23"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("E49").Value = "  -4.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.175"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07073"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.93%  "
